# Add edit login and password function
# Applies the data updates described by the diff: new tenant address, new
# payment date, new tenant (login/password holder) name, and updated
# individual/communal electricity meter readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Address of the premises (was: "г. Иркутск, ул. Лермонтова, д. 10, кв. 1")
$ws.Range("N2").Value = "г. Москва, ул. Пролетарская, д. 15, кв. 5"

# Date of the last received payment (was: "15.06.2021г.")
$ws.Range("O14").Value = "16.06.2021г."

# Consumer (tenant) full name (was: "Иванов Иван Иванович")
$ws.Range("D6").Value = "Астафьев Владимир Дмитриевич"

# Electricity meter readings: previous/current values
$ws.Range("E35").Value = 200.0
$ws.Range("G35").Value = 202.0
